$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fecha del Reporte: 45370 -> 45371
$ws.Range("C2").Value = 45371

# Ficha de Caracterización: "2824078" -> "2499992" (kept as text)
$ws.Range("C3").Value = "'2499992"

# Programa: "DISEÑO E INTEGRACION DE MULTIMEDIA" -> "PROGRAMA DE PRUEBA"
$ws.Range("C4").Value = "PROGRAMA DE PRUEBA"

# Fecha Inicio: 44760 -> 44669
$ws.Range("C6").Value = 44669

# Fecha Fin: 45124 -> 45490
$ws.Range("C7").Value = 45490
